# Applies the "Updated cryptos list on Fri Aug 23 19:47:54 UTC 2024 with GitHub Actions" data refresh.
# All cells in this sheet hold plain text (e.g. Price "63.261.72", Volume "  +5.07%  "), so
# numeric-looking Price values (column D) are written with a leading apostrophe - the classic
# Excel "force text" quote-prefix used when typing such values directly into a cell - so that
# Excel/COM doesn't silently reinterpret them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.261.72'
$ws.Range("E2").Value = '  +5.07%  '
$ws.Range("D3").Value = '''2.709.07'
$ws.Range("E3").Value = '  +4.23%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''581.37'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").Value = '''149.52'
$ws.Range("E6").Value = '  +4.67%  '
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("E8").Value = '  +1.40%  '
$ws.Range("D9").Value = '''2.739.03'
$ws.Range("E9").Value = '  +5.19%  '
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("E11").Value = '  +7.55%  '
$ws.Range("D12").Value = '''0.389'
$ws.Range("E12").Value = '  +4.83%  '
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").Value = '''3.187.88'
$ws.Range("E14").Value = '  +4.20%  '
$ws.Range("E15").Value = '  +8.54%  '
$ws.Range("D16").Value = '''63.171.04'
$ws.Range("E16").Value = '  +4.91%  '
$ws.Range("E17").Value = '  +7.34%  '
$ws.Range("D18").Value = '''2.728.18'
$ws.Range("E18").Value = '  +4.84%  '
$ws.Range("D19").Value = '''11.95'
$ws.Range("E19").Value = '  +5.53%  '
$ws.Range("E20").Value = '  +5.49%  '
$ws.Range("D21").Value = '''363.50'
$ws.Range("E21").Value = '  +5.22%  '
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").Value = '''0.996'
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").Value = '''0.534'
$ws.Range("E24").Value = '  +0.60%  '
$ws.Range("D25").Value = '''65.50'
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("E26").Value = '  +4.06%  '
$ws.Range("E27").Value = '  +6.96%  '
$ws.Range("D28").Value = '''0.995'
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.03'
$ws.Range("E29").Value = '  +6.83%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '''0.0₃0859'
$ws.Range("E30").Value = '  +7.78%  '
$ws.Range("E31").Value = '  +10.30%  '
$ws.Range("D32").Value = '''169.84'
$ws.Range("E32").Value = '  +1.82%  '
$ws.Range("E33").Value = '  +22.23%  '
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").Value = '''20.53'
$ws.Range("E35").Value = '  +5.76%  '
$ws.Range("E36").Value = '  +11.99%  '
$ws.Range("E37").Value = '  +9.00%  '
$ws.Range("E38").Value = '  +10.84%  '
$ws.Range("E39").Value = '  +19.16%  '
$ws.Range("D40").Value = '''351.87'
$ws.Range("E40").Value = '  +12.34%  '
$ws.Range("D41").Value = '''4.25'
$ws.Range("E41").Value = '  +9.86%  '
$ws.Range("D42").Value = '''39.32'
$ws.Range("E42").Value = '  +2.99%  '
$ws.Range("E43").Value = '  +13.19%  '
$ws.Range("D44").Value = '''21.61'
$ws.Range("E44").Value = '  +8.96%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = '''0.0591'
$ws.Range("E45").Value = '  +7.90%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''21.69'
$ws.Range("E46").Value = '  +9.22%  '
$ws.Range("D47").Value = '''139.04'
$ws.Range("E47").Value = '  +2.54%  '
$ws.Range("D48").Value = '''0.0259'
$ws.Range("E48").Value = '  +7.42%  '
$ws.Range("E49").Value = '  +5.69%  '
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("E51").Value = '  -0.41%  '
